# Update the "as_of_utc" timestamp column (AA2:AA26) on the "Главные" and
# "Линейные" sheets to reflect the new publish time.

$wb = $excel.ActiveWorkbook

$newTimestamp = "2025-11-13 09:48:49"
$sheetNames = @("Главные", "Линейные")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 2; $row -le 26; $row++) {
        $ws.Range("AA$row").Value = $newTimestamp
    }
}
